# Apply cryptos list update (prices / 1h volume %) for Wed Jul 10 18:36:23 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    # Force the cell to keep its original Text data type (matches the inlineStr cells
    # already in the sheet) instead of letting Excel auto-convert numeric-looking
    # strings (e.g. "522.70", "57.560.48") into real numbers, then restore the
    # cells style so no stray number-format style is left behind on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "57.560.48"
Set-TextCell $ws.Range("E2") "  -0.49%  "
Set-TextCell $ws.Range("D3") "3.099.40"
Set-TextCell $ws.Range("E3") "  +1.12%  "
Set-TextCell $ws.Range("E4") "  +0.01%  "
Set-TextCell $ws.Range("D5") "522.70"
Set-TextCell $ws.Range("E5") "  +0.50%  "
Set-TextCell $ws.Range("D6") "140.79"
Set-TextCell $ws.Range("E6") "  -1.02%  "
Set-TextCell $ws.Range("E7") "  +0.05%  "
Set-TextCell $ws.Range("D8") "3.099.07"
Set-TextCell $ws.Range("E8") "  +1.13%  "
Set-TextCell $ws.Range("E9") "  +0.24%  "
Set-TextCell $ws.Range("D10") "7.23"
Set-TextCell $ws.Range("E10") "  -0.72%  "
Set-TextCell $ws.Range("D11") "0.108"
Set-TextCell $ws.Range("E11") "  +0.77%  "
Set-TextCell $ws.Range("D12") "0.386"
Set-TextCell $ws.Range("E12") "  +2.48%  "
Set-TextCell $ws.Range("D13") "3.631.89"
Set-TextCell $ws.Range("E13") "  +1.24%  "
Set-TextCell $ws.Range("D14") "0.131"
Set-TextCell $ws.Range("E14") "  +1.35%  "
Set-TextCell $ws.Range("D15") "26.15"
Set-TextCell $ws.Range("E15") "  +1.40%  "
Set-TextCell $ws.Range("D16") "0.0000163"
Set-TextCell $ws.Range("E16") "  +0.22%  "
Set-TextCell $ws.Range("D17") "57.611.49"
Set-TextCell $ws.Range("E17") "  -0.48%  "
Set-TextCell $ws.Range("D18") "3.097.89"
Set-TextCell $ws.Range("E18") "  +1.43%  "
Set-TextCell $ws.Range("D19") "6.12"
Set-TextCell $ws.Range("E19") "  +0.72%  "
Set-TextCell $ws.Range("D20") "12.78"
Set-TextCell $ws.Range("E20") "  -0.71%  "
Set-TextCell $ws.Range("D21") "8.07"
Set-TextCell $ws.Range("E21") "  -1.05%  "
Set-TextCell $ws.Range("D22") "335.87"
Set-TextCell $ws.Range("E22") "  +1.55%  "
Set-TextCell $ws.Range("D23") "0.999"
Set-TextCell $ws.Range("E23") "  -0.01%  "
Set-TextCell $ws.Range("D24") "0.511"
Set-TextCell $ws.Range("E24") "  +2.44%  "
Set-TextCell $ws.Range("D25") "66.70"
Set-TextCell $ws.Range("E25") "  +1.60%  "
Set-TextCell $ws.Range("E26") "  -0.56%  "
Set-TextCell $ws.Range("E27") "  +0.32%  "
Set-TextCell $ws.Range("D28") "0.0₃0917"
Set-TextCell $ws.Range("E28") "  +1.88%  "
Set-TextCell $ws.Range("D29") "6.51"
Set-TextCell $ws.Range("E29") "  +2.01%  "
Set-TextCell $ws.Range("E30") "  -0.02%  "
Set-TextCell $ws.Range("D31") "7.21"
Set-TextCell $ws.Range("E31") "  +0.13%  "
Set-TextCell $ws.Range("D32") "1.86"
Set-TextCell $ws.Range("E32") "  +2.10%  "
Set-TextCell $ws.Range("D33") "20.99"
Set-TextCell $ws.Range("E33") "  +1.41%  "
Set-TextCell $ws.Range("E34") "  +0.36%  "
Set-TextCell $ws.Range("D35") "155.95"
Set-TextCell $ws.Range("E35") "  +0.70%  "
Set-TextCell $ws.Range("D36") "4.64"
Set-TextCell $ws.Range("E36") "  +2.71%  "
Set-TextCell $ws.Range("D37") "6.10"
Set-TextCell $ws.Range("E37") "  +2.41%  "
Set-TextCell $ws.Range("D38") "27.22"
Set-TextCell $ws.Range("E38") "  -0.40%  "
Set-TextCell $ws.Range("E39") "  +2.75%  "
Set-TextCell $ws.Range("D40") "0.0664"
Set-TextCell $ws.Range("E40") "  -1.53%  "
Set-TextCell $ws.Range("B41") "Mantle"
Set-TextCell $ws.Range("C41") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell $ws.Range("D41") "0.686"
Set-TextCell $ws.Range("E41") "  +5.10%  "
Set-TextCell $ws.Range("B42") "RenzoRestakedETH"
Set-TextCell $ws.Range("C42") "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
Set-TextCell $ws.Range("D42") "3.139.00"
Set-TextCell $ws.Range("E42") "  +1.14%  "
Set-TextCell $ws.Range("B43") "Filecoin"
Set-TextCell $ws.Range("C43") "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell $ws.Range("D43") "3.93"
Set-TextCell $ws.Range("E43") "  +0.20%  "
Set-TextCell $ws.Range("B44") "Stacks"
Set-TextCell $ws.Range("C44") "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell $ws.Range("D44") "1.51"
Set-TextCell $ws.Range("E44") "  +11.00%  "
Set-TextCell $ws.Range("D45") "36.79"
Set-TextCell $ws.Range("E45") "  +0.24%  "
Set-TextCell $ws.Range("D46") "1.00"
Set-TextCell $ws.Range("E46") "  -0.01%  "
Set-TextCell $ws.Range("D47") "2.311.20"
Set-TextCell $ws.Range("E47") "  +1.84%  "
Set-TextCell $ws.Range("E48") "  +0.79%  "
Set-TextCell $ws.Range("D49") "0.979"
Set-TextCell $ws.Range("E49") "  +5.72%  "
Set-TextCell $ws.Range("D50") "20.71"
Set-TextCell $ws.Range("E50") "  -0.60%  "
Set-TextCell $ws.Range("D51") "6.01"
Set-TextCell $ws.Range("E51") "  +2.09%  "
